$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates from the 2023-01-23 GitHub Actions symbol-list refresh.
# NumberFormat "@" + Style "Normal" keeps numeric-looking strings (prices,
# percentages) stored as literal text, matching the source data, while
# resetting the cell style back to the workbook default (no lingering
# Text-format style is left behind on the cell).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "305.37"
Set-TextValue "E2" "1.18%"
Set-TextValue "D3" "36.09"
Set-TextValue "D4" "5.097"
Set-TextValue "E4" "1.85%"
Set-TextValue "D5" "0.07857"
Set-TextValue "E5" "0.00%"
Set-TextValue "D6" "2.167"
Set-TextValue "E6" "-2.88%"
Set-TextValue "D7" "7.919"
Set-TextValue "E7" "-1.29%"
Set-TextValue "B8" "MXToken"
Set-TextValue "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.9218"
Set-TextValue "E8" "1.37%"
Set-TextValue "B9" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D9" "0.09682"
Set-TextValue "E9" "5.14%"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1867"
Set-TextValue "E10" "0.14%"
Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.08669"
Set-TextValue "E11" "2.02%"
Set-TextValue "B12" "BitrueCoin"
Set-TextValue "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.03484"
Set-TextValue "E12" "-0.94%"
Set-TextValue "B13" "BitMartToken"
Set-TextValue "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09928"
Set-TextValue "E13" "0.01%"
Set-TextValue "B14" "BitForexToken"
Set-TextValue "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001444"
Set-TextValue "E14" "-1.95%"
Set-TextValue "B15" "TigerCash"
Set-TextValue "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.005683"
Set-TextValue "E15" "0.30%"
Set-TextValue "B16" "LEO"
Set-TextValue "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.465"
Set-TextValue "E16" "-0.14%"
Set-TextValue "B17" "GateToken"
Set-TextValue "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D17" "4.097"
Set-TextValue "E17" "2.18%"
Set-TextValue "D18" "2.396"
Set-TextValue "E18" "16.28%"
Set-TextValue "D20" "0.1300"
Set-TextValue "E20" "-0.58%"
Set-TextValue "D21" "4.868"
Set-TextValue "E21" "1.40%"
Set-TextValue "E22" "-0.04%"
Set-TextValue "D23" "0.04554"
Set-TextValue "E23" "-2.09%"
Set-TextValue "D24" "0.005086"
Set-TextValue "E24" "14.32%"
Set-TextValue "E25" "0.38%"
Set-TextValue "D26" "0.0001401"
Set-TextValue "E26" "7.88%"
Set-TextValue "D27" "0.0004750"
Set-TextValue "E27" "0.14%"
Set-TextValue "D39" "0.01835"
Set-TextValue "E39" "4.23%"
Set-TextValue "D40" "0.04773"
Set-TextValue "E40" "1.18%"
Set-TextValue "D41" "0.007696"
Set-TextValue "E41" "-2.20%"
Set-TextValue "D42" "0.1398"
Set-TextValue "E42" "0.44%"
Set-TextValue "D43" "0.007740"
Set-TextValue "E43" "1.09%"
Set-TextValue "D44" "0.002231"
Set-TextValue "E44" "0.60%"
Set-TextValue "D45" "0.01133"
Set-TextValue "E45" "10.77%"
Set-TextValue "E46" "6.23%"
Set-TextValue "E47" "0.15%"
Set-TextValue "E48" "0.00%"
Set-TextValue "D49" "24.52"
Set-TextValue "E49" "182.82%"
Set-TextValue "E50" "-25.54%"
Set-TextValue "E51" "0.15%"
